# Sample Project / Main.xlsx — B11 (the "R40" rule-name cell) is changed
# to hold the text "1" instead, while keeping its existing cell style and
# storing the value as a shared string (not a number).
#
# A plain `Range.Value = "1"` would be auto-parsed as the number 1 (normal
# Excel "smart" typed-input behavior), and forcing text via NumberFormat
# "@" / a quote-prefix would mint a brand new cell style (bumping cellXfs)
# even after being reset back afterwards. Instead we write a formula that
# evaluates to the text "1", then convert that formula to a literal value
# via copy / paste-values, which keeps the original style untouched and
# stores a genuine text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
